$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage reordered rows into helper columns D:E, then copy back to A:B in one pass.
# This preserves original cell types/styles (e.g. text dates stay as shared-string text).
$ws.Range("A2:B2").Copy($ws.Range("D2:E2"))
$ws.Range("A3:B3").Copy($ws.Range("D3:E3"))
$ws.Range("A4:B4").Copy($ws.Range("D4:E4"))
$ws.Range("A5:B5").Copy($ws.Range("D5:E5"))
$ws.Range("A6:B6").Copy($ws.Range("D6:E6"))
$ws.Range("A7:B7").Copy($ws.Range("D7:E7"))
$ws.Range("A8:B8").Copy($ws.Range("D8:E8"))
$ws.Range("A9:B9").Copy($ws.Range("D9:E9"))
$ws.Range("A10:B10").Copy($ws.Range("D10:E10"))
$ws.Range("A11:B11").Copy($ws.Range("D11:E11"))
$ws.Range("A12:B12").Copy($ws.Range("D12:E12"))
$ws.Range("A13:B13").Copy($ws.Range("D13:E13"))
$ws.Range("A14:B14").Copy($ws.Range("D14:E14"))
$ws.Range("A15:B15").Copy($ws.Range("D15:E15"))
$ws.Range("A17:B17").Copy($ws.Range("D16:E16"))
$ws.Range("A16:B16").Copy($ws.Range("D17:E17"))
$ws.Range("A18:B18").Copy($ws.Range("D18:E18"))
$ws.Range("A19:B19").Copy($ws.Range("D19:E19"))
$ws.Range("A20:B20").Copy($ws.Range("D20:E20"))
$ws.Range("A22:B22").Copy($ws.Range("D21:E21"))
$ws.Range("A21:B21").Copy($ws.Range("D22:E22"))
$ws.Range("A23:B23").Copy($ws.Range("D23:E23"))
$ws.Range("A25:B25").Copy($ws.Range("D24:E24"))
$ws.Range("A24:B24").Copy($ws.Range("D25:E25"))
$ws.Range("A26:B26").Copy($ws.Range("D26:E26"))
$ws.Range("A27:B27").Copy($ws.Range("D27:E27"))
$ws.Range("A28:B28").Copy($ws.Range("D28:E28"))
$ws.Range("A30:B30").Copy($ws.Range("D29:E29"))
$ws.Range("A33:B33").Copy($ws.Range("D30:E30"))
$ws.Range("A32:B32").Copy($ws.Range("D31:E31"))
$ws.Range("A29:B29").Copy($ws.Range("D32:E32"))
$ws.Range("A31:B31").Copy($ws.Range("D33:E33"))
$ws.Range("A34:B34").Copy($ws.Range("D34:E34"))
$ws.Range("A35:B35").Copy($ws.Range("D35:E35"))
$ws.Range("A36:B36").Copy($ws.Range("D36:E36"))
$ws.Range("A37:B37").Copy($ws.Range("D37:E37"))
$ws.Range("A38:B38").Copy($ws.Range("D38:E38"))
$ws.Range("A39:B39").Copy($ws.Range("D39:E39"))
$ws.Range("A40:B40").Copy($ws.Range("D40:E40"))
$ws.Range("A43:B43").Copy($ws.Range("D41:E41"))
$ws.Range("A41:B41").Copy($ws.Range("D42:E42"))
$ws.Range("A42:B42").Copy($ws.Range("D43:E43"))
$ws.Range("A52:B52").Copy($ws.Range("D44:E44"))
$ws.Range("A46:B46").Copy($ws.Range("D45:E45"))
$ws.Range("A49:B49").Copy($ws.Range("D46:E46"))
$ws.Range("A51:B51").Copy($ws.Range("D47:E47"))
$ws.Range("A50:B50").Copy($ws.Range("D48:E48"))
$ws.Range("A48:B48").Copy($ws.Range("D49:E49"))
$ws.Range("A44:B44").Copy($ws.Range("D50:E50"))
$ws.Range("A45:B45").Copy($ws.Range("D51:E51"))
$ws.Range("A47:B47").Copy($ws.Range("D52:E52"))
$ws.Range("A53:B53").Copy($ws.Range("D53:E53"))
$ws.Range("A54:B54").Copy($ws.Range("D54:E54"))
$ws.Range("A55:B55").Copy($ws.Range("D55:E55"))
$ws.Range("A56:B56").Copy($ws.Range("D56:E56"))
$ws.Range("A57:B57").Copy($ws.Range("D57:E57"))
$ws.Range("A58:B58").Copy($ws.Range("D58:E58"))
$ws.Range("A60:B60").Copy($ws.Range("D59:E59"))
$ws.Range("A59:B59").Copy($ws.Range("D60:E60"))
$ws.Range("A62:B62").Copy($ws.Range("D61:E61"))
$ws.Range("A61:B61").Copy($ws.Range("D62:E62"))
$ws.Range("A63:B63").Copy($ws.Range("D63:E63"))
$ws.Range("A64:B64").Copy($ws.Range("D64:E64"))
$ws.Range("A65:B65").Copy($ws.Range("D65:E65"))
$ws.Range("A66:B66").Copy($ws.Range("D66:E66"))
$ws.Range("A68:B68").Copy($ws.Range("D67:E67"))
$ws.Range("A67:B67").Copy($ws.Range("D68:E68"))
$ws.Range("A69:B69").Copy($ws.Range("D69:E69"))
$ws.Range("A70:B70").Copy($ws.Range("D70:E70"))
$ws.Range("A72:B72").Copy($ws.Range("D71:E71"))
$ws.Range("A71:B71").Copy($ws.Range("D72:E72"))
$ws.Range("A74:B74").Copy($ws.Range("D73:E73"))
$ws.Range("A75:B75").Copy($ws.Range("D74:E74"))
$ws.Range("A73:B73").Copy($ws.Range("D75:E75"))
$ws.Range("A77:B77").Copy($ws.Range("D76:E76"))
$ws.Range("A78:B78").Copy($ws.Range("D77:E77"))
$ws.Range("A76:B76").Copy($ws.Range("D78:E78"))
$ws.Range("A80:B80").Copy($ws.Range("D79:E79"))
$ws.Range("A79:B79").Copy($ws.Range("D80:E80"))
$ws.Range("A83:B83").Copy($ws.Range("D81:E81"))
$ws.Range("A81:B81").Copy($ws.Range("D82:E82"))
$ws.Range("A84:B84").Copy($ws.Range("D83:E83"))
$ws.Range("A82:B82").Copy($ws.Range("D84:E84"))
$ws.Range("A85:B85").Copy($ws.Range("D85:E85"))
$ws.Range("A86:B86").Copy($ws.Range("D86:E86"))
$ws.Range("A87:B87").Copy($ws.Range("D87:E87"))
$ws.Range("A89:B89").Copy($ws.Range("D88:E88"))
$ws.Range("A88:B88").Copy($ws.Range("D89:E89"))
$ws.Range("A91:B91").Copy($ws.Range("D90:E90"))
$ws.Range("A90:B90").Copy($ws.Range("D91:E91"))
$ws.Range("A92:B92").Copy($ws.Range("D92:E92"))
$ws.Range("A93:B93").Copy($ws.Range("D93:E93"))
$ws.Range("A94:B94").Copy($ws.Range("D94:E94"))
$ws.Range("A97:B97").Copy($ws.Range("D95:E95"))
$ws.Range("A96:B96").Copy($ws.Range("D96:E96"))
$ws.Range("A99:B99").Copy($ws.Range("D97:E97"))
$ws.Range("A95:B95").Copy($ws.Range("D98:E98"))
$ws.Range("A98:B98").Copy($ws.Range("D99:E99"))
$ws.Range("A101:B101").Copy($ws.Range("D100:E100"))
$ws.Range("A100:B100").Copy($ws.Range("D101:E101"))
$ws.Range("A104:B104").Copy($ws.Range("D102:E102"))
$ws.Range("A106:B106").Copy($ws.Range("D103:E103"))
$ws.Range("A107:B107").Copy($ws.Range("D104:E104"))
$ws.Range("A111:B111").Copy($ws.Range("D105:E105"))
$ws.Range("A103:B103").Copy($ws.Range("D106:E106"))
$ws.Range("A112:B112").Copy($ws.Range("D107:E107"))
$ws.Range("A102:B102").Copy($ws.Range("D108:E108"))
$ws.Range("A108:B108").Copy($ws.Range("D109:E109"))
$ws.Range("A105:B105").Copy($ws.Range("D110:E110"))
$ws.Range("A110:B110").Copy($ws.Range("D111:E111"))
$ws.Range("A109:B109").Copy($ws.Range("D112:E112"))
$ws.Range("A113:B113").Copy($ws.Range("D113:E113"))
$ws.Range("A114:B114").Copy($ws.Range("D114:E114"))
$ws.Range("A115:B115").Copy($ws.Range("D115:E115"))
$ws.Range("A119:B119").Copy($ws.Range("D116:E116"))
$ws.Range("A118:B118").Copy($ws.Range("D117:E117"))
$ws.Range("A117:B117").Copy($ws.Range("D118:E118"))
$ws.Range("A120:B120").Copy($ws.Range("D119:E119"))
$ws.Range("A116:B116").Copy($ws.Range("D120:E120"))
$ws.Range("A121:B121").Copy($ws.Range("D121:E121"))
$ws.Range("A122:B122").Copy($ws.Range("D122:E122"))
$ws.Range("A123:B123").Copy($ws.Range("D123:E123"))
$ws.Range("A124:B124").Copy($ws.Range("D124:E124"))
$ws.Range("A125:B125").Copy($ws.Range("D125:E125"))
$ws.Range("A127:B127").Copy($ws.Range("D126:E126"))
$ws.Range("A126:B126").Copy($ws.Range("D127:E127"))
$ws.Range("A128:B128").Copy($ws.Range("D128:E128"))
$ws.Range("A129:B129").Copy($ws.Range("D129:E129"))
$ws.Range("A130:B130").Copy($ws.Range("D130:E130"))
$ws.Range("A131:B131").Copy($ws.Range("D131:E131"))
$ws.Range("A132:B132").Copy($ws.Range("D132:E132"))
$ws.Range("A133:B133").Copy($ws.Range("D133:E133"))
$ws.Range("A135:B135").Copy($ws.Range("D134:E134"))
$ws.Range("A134:B134").Copy($ws.Range("D135:E135"))
$ws.Range("A137:B137").Copy($ws.Range("D136:E136"))
$ws.Range("A136:B136").Copy($ws.Range("D137:E137"))
$ws.Range("A138:B138").Copy($ws.Range("D138:E138"))
$ws.Range("A140:B140").Copy($ws.Range("D139:E139"))
$ws.Range("A139:B139").Copy($ws.Range("D140:E140"))
$ws.Range("A142:B142").Copy($ws.Range("D141:E141"))
$ws.Range("A143:B143").Copy($ws.Range("D142:E142"))
$ws.Range("A144:B144").Copy($ws.Range("D143:E143"))
$ws.Range("A141:B141").Copy($ws.Range("D144:E144"))
$ws.Range("A147:B147").Copy($ws.Range("D145:E145"))
$ws.Range("A148:B148").Copy($ws.Range("D146:E146"))
$ws.Range("A145:B145").Copy($ws.Range("D147:E147"))
$ws.Range("A146:B146").Copy($ws.Range("D148:E148"))
$ws.Range("A149:B149").Copy($ws.Range("D149:E149"))
$ws.Range("A155:B155").Copy($ws.Range("D150:E150"))
$ws.Range("A150:B150").Copy($ws.Range("D151:E151"))
$ws.Range("A153:B153").Copy($ws.Range("D152:E152"))
$ws.Range("A151:B151").Copy($ws.Range("D153:E153"))
$ws.Range("A154:B154").Copy($ws.Range("D154:E154"))
$ws.Range("A152:B152").Copy($ws.Range("D155:E155"))
$ws.Range("A158:B158").Copy($ws.Range("D156:E156"))
$ws.Range("A156:B156").Copy($ws.Range("D157:E157"))
$ws.Range("A157:B157").Copy($ws.Range("D158:E158"))
$ws.Range("A161:B161").Copy($ws.Range("D159:E159"))
$ws.Range("A159:B159").Copy($ws.Range("D160:E160"))
$ws.Range("A162:B162").Copy($ws.Range("D161:E161"))
$ws.Range("A160:B160").Copy($ws.Range("D162:E162"))
$ws.Range("A164:B164").Copy($ws.Range("D163:E163"))
$ws.Range("A163:B163").Copy($ws.Range("D164:E164"))
$ws.Range("A165:B165").Copy($ws.Range("D165:E165"))
$ws.Range("A168:B168").Copy($ws.Range("D166:E166"))
$ws.Range("A166:B166").Copy($ws.Range("D167:E167"))
$ws.Range("A167:B167").Copy($ws.Range("D168:E168"))
$ws.Range("A171:B171").Copy($ws.Range("D169:E169"))
$ws.Range("A169:B169").Copy($ws.Range("D170:E170"))
$ws.Range("A172:B172").Copy($ws.Range("D171:E171"))
$ws.Range("A170:B170").Copy($ws.Range("D172:E172"))
$ws.Range("A174:B174").Copy($ws.Range("D173:E173"))
$ws.Range("A173:B173").Copy($ws.Range("D174:E174"))
$ws.Range("A175:B175").Copy($ws.Range("D175:E175"))
$ws.Range("A176:B176").Copy($ws.Range("D176:E176"))
$ws.Range("A177:B177").Copy($ws.Range("D177:E177"))
$ws.Range("A178:B178").Copy($ws.Range("D178:E178"))
$ws.Range("A180:B180").Copy($ws.Range("D179:E179"))
$ws.Range("A179:B179").Copy($ws.Range("D180:E180"))
$ws.Range("A181:B181").Copy($ws.Range("D181:E181"))
$ws.Range("A182:B182").Copy($ws.Range("D182:E182"))
$ws.Range("A183:B183").Copy($ws.Range("D183:E183"))
$ws.Range("A184:B184").Copy($ws.Range("D184:E184"))
$ws.Range("A185:B185").Copy($ws.Range("D185:E185"))
$ws.Range("A186:B186").Copy($ws.Range("D186:E186"))

# Copy the staged block back over the original data range in one operation.
$ws.Range("D2:E186").Copy($ws.Range("A2:B186"))

# Clear the staging area so it does not appear in the used range.
$ws.Range("D2:E186").Clear()
